$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new rows of data (date in column A, count in column B),
# matching the date-formatted style already used by the rows above.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A17").PasteSpecial(-4122) | Out-Null
$ws.Range("A17").Value = 45877
$ws.Range("B17").Value = 9

$ws.Range("A2").Copy() | Out-Null
$ws.Range("A18").PasteSpecial(-4122) | Out-Null
$ws.Range("A18").Value = 45894
$ws.Range("B18").Value = 56

# Auto-fit column A to its (now longer) content.
$ws.Columns("A:A").AutoFit() | Out-Null

# Leave the active selection on G17, like the saved workbook shows.
$ws.Range("G17").Select() | Out-Null
